# Applies the coinranking.com market-data refresh for Thu Dec 21 2023 GitHub
# Actions run: updated Price (col D) / Volume(1h) (col E) figures, and rows 46-47
# swap places because Aave overtook Cronos in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "43.885.29"
$ws.Range("E2").Value = "  -0.64%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.230.35"
$ws.Range("E3").Value = "  -0.98%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.03%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'274.29"
$ws.Range("E5").Value = "  +6.23%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "'87.89"
$ws.Range("E6").Value = "  +8.66%  "

# --- Row 7: XRP ---
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  -1.22%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.07%  "

# --- Row 10: Avalanche ---
$ws.Range("D10").Value = "'45.25"
$ws.Range("E10").Value = "  +3.84%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  -1.52%  "

# --- Row 12: Polkadot ---
$ws.Range("D12").Value = "'7.68"
$ws.Range("E12").Value = "  +7.93%  "

# --- Row 13: TRON ---
$ws.Range("E13").Value = "  +1.06%  "

# --- Row 14: WrappedliquidstakedEther2.0 ---
$ws.Range("D14").Value = "2.563.29"
$ws.Range("E14").Value = "  -0.86%  "

# --- Row 15: Chainlink ---
$ws.Range("D15").Value = "'14.98"
$ws.Range("E15").Value = "  +1.03%  "

# --- Row 16: WrappedEther ---
$ws.Range("D16").Value = "2.252.23"
$ws.Range("E16").Value = "  +1.08%  "

# --- Row 17: Polygon ---
$ws.Range("D17").Value = "'0.791"
$ws.Range("E17").Value = "  -0.85%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "43.775.61"
$ws.Range("E18").Value = "  -0.66%  "

# --- Row 19: ShibaInu ---
$ws.Range("E19").Value = "  -1.31%  "

# --- Row 20: Litecoin ---
$ws.Range("D20").Value = "'70.24"
$ws.Range("E20").Value = "  -2.09%  "

# --- Row 21: Uniswap ---
$ws.Range("D21").Value = "'5.97"
$ws.Range("E21").Value = "  -1.85%  "

# --- Row 22: ImmutableX ---
$ws.Range("E22").Value = "  -0.09%  "

# --- Row 23: BitcoinCash ---
$ws.Range("D23").Value = "'232.61"
$ws.Range("E23").Value = "  -1.20%  "

# --- Row 24: InternetComputer(DFINITY) ---
$ws.Range("D24").Value = "'8.75"
$ws.Range("E24").Value = "  -7.83%  "

# --- Row 25: PancakeSwap ---
$ws.Range("E25").Value = "  +14.65%  "

# --- Row 26: Dai ---
$ws.Range("E26").Value = "  -0.06%  "

# --- Row 27: Cosmos ---
$ws.Range("E27").Value = "  -0.93%  "

# --- Row 28: WEMIXToken ---
$ws.Range("E28").Value = "  +3.71%  "

# --- Row 29: Toncoin ---
$ws.Range("E29").Value = "  +4.24%  "

# --- Row 30: InjectiveProtocol ---
$ws.Range("D30").Value = "'39.12"
$ws.Range("E30").Value = "  -4.82%  "

# --- Row 31: Monero ---
$ws.Range("D31").Value = "'172.69"
$ws.Range("E31").Value = "  -0.42%  "

# --- Row 32: Hedera ---
$ws.Range("D32").Value = "'0.0902"
$ws.Range("E32").Value = "  +2.70%  "

# --- Row 33: EthereumClassic ---
$ws.Range("D33").Value = "'20.77"
$ws.Range("E33").Value = "  +0.37%  "

# --- Row 35: Stellar ---
$ws.Range("E35").Value = "  -0.08%  "

# --- Row 36: Kaspa ---
$ws.Range("E36").Value = "  -3.67%  "

# --- Row 37: VeChain ---
$ws.Range("D37").Value = "'0.0353"
$ws.Range("E37").Value = "  -4.21%  "

# --- Row 38: RenderToken ---
$ws.Range("D38").Value = "'4.26"
$ws.Range("E38").Value = "  -5.97%  "

# --- Row 39: NEARProtocol ---
$ws.Range("D39").Value = "'3.45"
$ws.Range("E39").Value = "  +17.56%  "

# --- Row 40: Celestia ---
$ws.Range("D40").Value = "'12.44"
$ws.Range("E40").Value = "  -5.09%  "

# --- Row 41: LidoDAOToken ---
$ws.Range("D41").Value = "'2.16"
$ws.Range("E41").Value = "  +0.29%  "

# --- Row 42: MultiversX ---
$ws.Range("D42").Value = "'63.91"
$ws.Range("E42").Value = "  +0.92%  "

# --- Row 43: Algorand ---
$ws.Range("D43").Value = "'0.209"
$ws.Range("E43").Value = "  +1.70%  "

# --- Row 44: THORChain ---
$ws.Range("D44").Value = "'5.39"
$ws.Range("E44").Value = "  -3.05%  "

# --- Row 45: FraxShare ---
$ws.Range("D45").Value = "'8.50"
$ws.Range("E45").Value = "  -1.02%  "

# --- Rows 46 & 47: Aave overtakes Cronos, so the two rows swap contents ---
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'100.45"
$ws.Range("E46").Value = "  -3.73%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0982"
$ws.Range("E47").Value = "  -1.47%  "

# --- Row 48: TrustWalletToken ---
$ws.Range("D48").Value = "'1.20"
$ws.Range("E48").Value = "  +3.51%  "

# --- Row 49: ARBITRUM ---
$ws.Range("E49").Value = "  +0.03%  "

# --- Row 50: Stacks ---
$ws.Range("D50").Value = "'1.50"
$ws.Range("E50").Value = "  -2.64%  "

# --- Row 51: WOONetwork ---
$ws.Range("D51").Value = "'0.429"
$ws.Range("E51").Value = "  -7.60%  "
